# Add a new "ListBullet" paragraph listing the responsible professor right
# after the "Docente(s) Responsável(eis) " heading paragraph.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd("`r", "`a")
    if ($text -eq "Docente(s) Responsável(eis) ") {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not find the 'Docente(s) Responsável(eis)' paragraph"
}

# Insert a brand-new paragraph right after the heading, then style/fill it.
$target.Range.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Range.Text = "4893449 - Débora Souza Alvim"
$newPara.Style = "ListBullet"
